$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new bullet "Составлена диаграмма активности;" right
#    before the paragraph that begins "Частично реализованы классы
#    в решении ..." (same list: pStyle a3 / numId 11).
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Частично реализованы классы в решении")) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'Частично реализованы классы в решении' paragraph"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "Составлена диаграмма активности;"

# ------------------------------------------------------------------
# 2) Fix the split word "рели" + (bookmark) + "зный" -> "релизный",
#    and relocate the (hidden) _GoBack bookmark into its own new,
#    non-list paragraph right after "Составлен отчетный документ по
#    ролям (релизный документ)." and before the trailing empty
#    paragraph.
# ------------------------------------------------------------------

# Remove the bookmark from its old position (leaves surrounding runs
# untouched).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Join "рели" + "зный" into a single word "релизный". Delete the
# second fragment first so the search for "рели" afterwards can't
# accidentally match the "зный" tail that "релизный" would contain.
$d.Content.Find.Execute("зный", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("рели", $true, $false, $false, $false, $false, $true, 1, $false, "релизный", 2) | Out-Null

# Find the trailing empty paragraph (last paragraph of the document)
# and insert a fresh paragraph before it with the same (non-list)
# formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$bmParaIndex = $d.Paragraphs.Count - 1
$bmPara = $d.Paragraphs.Item($bmParaIndex)

# Use a throwaway placeholder character so we can anchor a
# non-collapsed Range (collapsed ranges right at a paragraph mark are
# not placed correctly by Bookmarks.Add), then add the bookmark
# around it and remove the placeholder again.
$insPos = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)
$insPos.InsertBefore("X")

$bmPara2 = $d.Paragraphs.Item($bmParaIndex)
$bmRange = $d.Range($bmPara2.Range.Start, $bmPara2.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmPara2.Range.Start, $bmPara2.Range.Start + 1)
$placeholder.Text = ""

Write-Output "done"
